# Add the team's season record (Wins / Losses / Ties) as three new
# trailing columns (AD, AE, AF) on the single worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 55

$wins = 101
$losses = 61
$ties = 0

# Header row (row 1) — new column headers.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header formatting (bold, thin border box, centered)
# by copying the format from the neighboring header cell AC1, which
# already carries that style.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows (2..55) — same season record repeated for every player.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # AD
    $ws.Cells.Item($r, 31).Value = $losses  # AE
    $ws.Cells.Item($r, 32).Value = $ties    # AF
}
